# ldlc_suivi_smartphones: append a new price-history snapshot column.
#
# The sheet keeps one column per scrape timestamp (A.. up through GB),
# followed by two fixed trailing columns: "nom" (product name) and
# "url_produit" (product url). Each new run inserts a fresh timestamp
# column right before "nom", shifting "nom"/"url_produit" one column to
# the right, and seeds the new column with the most recent known price
# for every product row (i.e. a straight copy of the previous last
# price column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "nom" currently lives in column GC (185th column). Inserting here
# pushes GC -> GD (nom) and GD -> GE (url_produit), and grows the used
# range from A1:GD210 to A1:GE210.
$nomCol = 185

$ws.Columns.Item($nomCol).Insert()

# Header for the freshly inserted timestamp column.
$ws.Cells.Item(1, $nomCol).Value = "2026-02-05 13:58:32"

# Seed every product row's new snapshot with the last recorded price
# (the column immediately to the left, which used to be the final price
# column before "nom" got shifted away). Rows with no recorded price yet
# simply copy the same blank.
$lastPriceCol = $nomCol - 1
$lastRow = 210

for ($row = 2; $row -le $lastRow; $row++) {
    $previousPrice = $ws.Cells.Item($row, $lastPriceCol).Value()
    $ws.Cells.Item($row, $nomCol).Value = $previousPrice
}
